# jbkowp/src/test/resources/OwpTestData.xlsx
# "added separate report and screenshot code"
#
# The D4 test-data cell ("123456") is no longer needed now that the
# report/screenshot code path is separated out, so its content is
# cleared (keeping its existing cell formatting), the active selection
# on the Login sheet moves back to D4, and the workbook window is
# resized slightly wider.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Remove the now-unused "123456" value from D4, keeping the cell's
# existing style/formatting intact.
$ws.Range("D4").ClearContents()

# Move the sheet's active selection to D4.
$ws.Range("D4").Select()

# Widen the saved workbook window.
$win = $wb.Windows.Item(1)
$win.Left = 0
$win.Top = 2520
$win.Width = 15360
$win.Height = 6330
